$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 43 ----
$ws.Range("A43").Value = 42615
$ws.Range("A43").NumberFormat = "d-mmm-yy"

$ws.Range("B43").Value = "CHPC"
$ws.Range("B43").NumberFormat = "d-mmm-yy"

$ws.Range("C43").Value = "Sim+analysis"

$ws.Range("D43").Value = 31
$ws.Range("E43").Value = 1254
$ws.Range("F43").Value = 10000

$ws.Range("G43").Formula = "=F43*E43"

$ws.Range("H43").Value = 0.66
$ws.Range("H43").HorizontalAlignment = -4152

$ws.Range("I43").Value = 1
$ws.Range("I43").HorizontalAlignment = -4152

# ---- Row 44 ----
$ws.Range("A44").Value = 42615
$ws.Range("A44").NumberFormat = "d-mmm-yy"

$ws.Range("B44").Value = "CHPC"
$ws.Range("B44").NumberFormat = "d-mmm-yy"

$ws.Range("C44").Value = "Sim+analysis"

$ws.Range("D44").Value = 31
$ws.Range("E44").Value = 1254
$ws.Range("F44").Value = 30000

$ws.Range("G44").Formula = "=F44*E44"

$ws.Range("H44").Value = 2
$ws.Range("H44").HorizontalAlignment = -4152

$ws.Range("I44").Value = 1
$ws.Range("I44").HorizontalAlignment = -4152

# ---- Selection as left by the author ----
$ws.Range("G44").Select()
